$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jan-2024")

# --- Add a new attendance column for 7 Jan 2024 (column I) ---
# Clone the existing date-column (G) formatting into I, then the blank
# trailing columns (J) formatting into the newly added K:M columns, so the
# new cells pick up the same styles already used on the sheet.
$ws.Range("G1:G8").Copy()
$ws.Range("I1:I8").PasteSpecial(-4122)
$ws.Range("J1:J8").Copy()
$ws.Range("K1:M8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header date for the new column: 7 Jan 2024
$ws.Range("I1").Value = 45298

# Attendance for 7 Jan 2024
$ws.Range("I2").Value = "Present"
$ws.Range("I3").Value = "Reason"
$ws.Range("I4").Value = "Absent"
$ws.Range("I5").Value = "Reason"
$ws.Range("I6").Value = "Present"
$ws.Range("I7").Value = "Present"
$ws.Range("I8").Value = "Absent"

# Reasons noted as threaded comments
$ws.Range("I3").AddCommentThreaded("Not feeling well reason cold and fever`n")
$ws.Range("I5").AddCommentThreaded("Medical emergency")

# Extend the Present/Absent/Reason dropdown validation to cover the new columns
$ws.Range("C2:M8").Validation.Delete()
$ws.Range("C2:M8").Validation.Add(3, 1, 1, '"Present, Absent,Reason"')

# Match the author's final selection
$ws.Range("I8").Select()
